$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 33200
$ws.Range("C3").Value = 61600
$ws.Range("D3").Value = 74100
$ws.Range("E3").Value = 68600
$ws.Range("B4").Value = 136
$ws.Range("C4").Value = 252
$ws.Range("D4").Value = 304
$ws.Range("E4").Value = 281
$ws.Range("B5").Value = 29.93895
$ws.Range("C5").Value = 31.87626
$ws.Range("D5").Value = 40.20043
$ws.Range("E5").Value = 56.1516
$ws.Range("B11").Value = 3047
$ws.Range("C11").Value = 21100
$ws.Range("D11").Value = 8062
$ws.Range("E11").Value = 33600
$ws.Range("B12").Value = 399
$ws.Range("C12").Value = 2767
$ws.Range("D12").Value = 1057
$ws.Range("E12").Value = 4401
$ws.Range("B13").Value = 327.04
$ws.Range("C13").Value = 93.88
$ws.Range("D13").Value = 363.11
$ws.Range("E13").Value = 116.82
$ws.Range("B19").Value = 356000
$ws.Range("C19").Value = 493000
$ws.Range("D19").Value = 578000
$ws.Range("E19").Value = 615000
$ws.Range("B20").Value = 1459
$ws.Range("C20").Value = 2018
$ws.Range("D20").Value = 2369
$ws.Range("E20").Value = 2521
$ws.Range("B21").Value = 1.24245
$ws.Range("C21").Value = 1.23562
$ws.Range("D21").Value = 1.30459
$ws.Range("E21").Value = 1.28014
$ws.Range("B27").Value = 14200
$ws.Range("C27").Value = 19300
$ws.Range("D27").Value = 21000
$ws.Range("E27").Value = 22500
$ws.Range("B28").Value = 1864
$ws.Range("C28").Value = 2532
$ws.Range("D28").Value = 2758
$ws.Range("E28").Value = 2950
$ws.Range("B29").Value = 23.7
$ws.Range("C29").Value = 26.14
$ws.Range("D29").Value = 28.42
$ws.Range("E29").Value = 31.46
$ws.Range("B35").Value = 529000
$ws.Range("C35").Value = 790000
$ws.Range("D35").Value = 902000
$ws.Range("E35").Value = 971000
$ws.Range("B36").Value = 2165
$ws.Range("C36").Value = 3234
$ws.Range("D36").Value = 3694
$ws.Range("E36").Value = 3977
$ws.Range("B37").Value = 1.68315
$ws.Range("C37").Value = 2.42806
$ws.Range("D37").Value = 3.19291
$ws.Range("E37").Value = 3.97811
$ws.Range("B43").Value = 20900
$ws.Range("C43").Value = 28400
$ws.Range("D43").Value = 31700
$ws.Range("E43").Value = 34700
$ws.Range("B44").Value = 2739
$ws.Range("C44").Value = 3728
$ws.Range("D44").Value = 4151
$ws.Range("E44").Value = 4550
$ws.Range("B45").Value = 47.28
$ws.Range("C45").Value = 69.7
$ws.Range("D45").Value = 91.81999999999999
$ws.Range("E45").Value = 112.23
$ws.Range("B51").Value = 280000
$ws.Range("C51").Value = 428000
$ws.Range("D51").Value = 445000
$ws.Range("E51").Value = 478000
$ws.Range("B52").Value = 1147
$ws.Range("C52").Value = 1754
$ws.Range("D52").Value = 1822
$ws.Range("E52").Value = 1959
$ws.Range("B53").Value = 1.39846
$ws.Range("C53").Value = 1.18719
$ws.Range("D53").Value = 1.18228
$ws.Range("E53").Value = 1.21185
$ws.Range("B59").Value = 11800
$ws.Range("C59").Value = 13700
$ws.Range("D59").Value = 14500
$ws.Range("E59").Value = 15400
$ws.Range("B60").Value = 1543
$ws.Range("C60").Value = 1790
$ws.Range("D60").Value = 1899
$ws.Range("E60").Value = 2018
$ws.Range("B61").Value = 23.74469
$ws.Range("C61").Value = 25.10921
$ws.Range("D61").Value = 27.41
$ws.Range("E61").Value = 30.15
$ws.Range("B67").Value = 32300
$ws.Range("C67").Value = 61700
$ws.Range("D67").Value = 52600
$ws.Range("E67").Value = 73300
$ws.Range("B68").Value = 132
$ws.Range("C68").Value = 253
$ws.Range("D68").Value = 216
$ws.Range("E68").Value = 300
$ws.Range("B69").Value = 30.8725
$ws.Range("C69").Value = 31.8884
$ws.Range("D69").Value = 56.5564
$ws.Range("E69").Value = 53.89756
$ws.Range("B75").Value = 11000
$ws.Range("C75").Value = 7816
$ws.Range("D75").Value = 28200
$ws.Range("E75").Value = 30600
$ws.Range("B76").Value = 1443
$ws.Range("C76").Value = 1025
$ws.Range("D76").Value = 3694
$ws.Range("E76").Value = 4006
$ws.Range("B77").Value = 90.59999999999999
$ws.Range("C77").Value = 254.41
$ws.Range("D77").Value = 105.33
$ws.Range("E77").Value = 127.87
$ws.Range("B83").Value = 315000
$ws.Range("C83").Value = 546000
$ws.Range("D83").Value = 626000
$ws.Range("E83").Value = 662000
$ws.Range("B84").Value = 1291
$ws.Range("C84").Value = 2237
$ws.Range("D84").Value = 2565
$ws.Range("E84").Value = 2711
$ws.Range("B85").Value = 1.47591
$ws.Range("C85").Value = 1.27465
$ws.Range("D85").Value = 1.22348
$ws.Range("E85").Value = 1.22044
$ws.Range("B91").Value = 16000
$ws.Range("C91").Value = 18300
$ws.Range("D91").Value = 21300
$ws.Range("E91").Value = 22600
$ws.Range("B92").Value = 2097
$ws.Range("C92").Value = 2397
$ws.Range("D92").Value = 2796
$ws.Range("E92").Value = 2966
$ws.Range("B93").Value = 22.11371
$ws.Range("C93").Value = 27.55
$ws.Range("D93").Value = 28.77
$ws.Range("E93").Value = 30.96
$ws.Range("B99").Value = 575000
$ws.Range("C99").Value = 771000
$ws.Range("D99").Value = 870000
$ws.Range("E99").Value = 964000
$ws.Range("B100").Value = 2355
$ws.Range("C100").Value = 3158
$ws.Range("D100").Value = 3563
$ws.Range("E100").Value = 3948
$ws.Range("B101").Value = 1.6196
$ws.Range("C101").Value = 2.44405
$ws.Range("D101").Value = 3.27233
$ws.Range("E101").Value = 3.99944
$ws.Range("B107").Value = 20900
$ws.Range("C107").Value = 26900
$ws.Range("D107").Value = 32300
$ws.Range("E107").Value = 35000
$ws.Range("B108").Value = 2739
$ws.Range("C108").Value = 3532
$ws.Range("D108").Value = 4238
$ws.Range("E108").Value = 4589
$ws.Range("B109").Value = 47.59
$ws.Range("C109").Value = 71.41
$ws.Range("D109").Value = 91.63
$ws.Range("E109").Value = 112.64
$ws.Range("B115").Value = 282000
$ws.Range("C115").Value = 420000
$ws.Range("D115").Value = 468000
$ws.Range("E115").Value = 489000
$ws.Range("B116").Value = 1157
$ws.Range("C116").Value = 1721
$ws.Range("D116").Value = 1917
$ws.Range("E116").Value = 2003
$ws.Range("B117").Value = 1.29361
$ws.Range("C117").Value = 1.12045
$ws.Range("D117").Value = 1.10408
$ws.Range("E117").Value = 1.11008
$ws.Range("B123").Value = 10800
$ws.Range("C123").Value = 13900
$ws.Range("D123").Value = 14900
$ws.Range("E123").Value = 15300
$ws.Range("B124").Value = 1413
$ws.Range("C124").Value = 1826
$ws.Range("D124").Value = 1955
$ws.Range("E124").Value = 2011
$ws.Range("B125").Value = 23.53156
$ws.Range("C125").Value = 28.30915
$ws.Range("D125").Value = 27.83776
$ws.Range("E125").Value = 31.47757
